# Add Individual Student Generate/Send Admit Card functionality
# Adds a "Result" worksheet (placed after "Sheet1") that mirrors the
# HUTOPSIds column from Sheet1 and appends a status column reporting
# whether each student's result was updated or the record was not found.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Number of populated rows (header + data) in Sheet1 column A.
$lastRow = $ws1.UsedRange.Rows.Count

# Create the new "Result" sheet immediately after Sheet1.
$resultSheet = $wb.Worksheets.Add($null, $ws1)
$resultSheet.Name = "Result"

# Status values for column B, aligned row-by-row with column A.
$statuses = @(
    "Record not found ",
    "Result Updated Successfully",
    "Result Updated Successfully",
    "Result Updated Successfully",
    "Result Updated Successfully",
    "Result Updated Successfully",
    "Result Updated Successfully",
    "Result Updated Successfully",
    "Result Updated Successfully",
    "Record not found ",
    "Result Updated Successfully"
)

for ($r = 1; $r -le $lastRow; $r++) {
    $idValue = $ws1.Cells.Item($r, 1).Text
    $resultSheet.Cells.Item($r, 1).Value = $idValue
    $resultSheet.Cells.Item($r, 2).Value = $statuses[$r - 1]
}

$wb.Save()
